$wb = $excel.ActiveWorkbook

# --- Commands sheet: Phidget HUB PWM "channel" -> "port" terminology update ---
$wsCommands = $wb.Worksheets.Item("Commands")

$wsCommands.Range("B66").Value = "outhub(<port>,<value>[,<sn>])"
$wsCommands.Range("C66").Value = "PHIDGET HUB PWM Output ON port <port> to  <value> in [0-100]"
$wsCommands.Range("B67").Value = "togglehub(<port>[,<sn>])"
$wsCommands.Range("C67").Value = "PHIDGET HUB PWM Output: toggles <port>"
$wsCommands.Range("B68").Value = "pulsehub(<port>,<millis>[,<sn>])"
$wsCommands.Range("C68").Value = "PHIDGET HUB PWM Output:  turn <port> ON for <millis> milliseconds"

# Narrow column A on the Commands sheet (was 112.84 chars wide, now much narrower)
$wsCommands.Columns.Item(1).ColumnWidth = 25.83

# --- Window / selection state ---
# The Commands sheet becomes the active tab, with C66 selected and the view
# scrolled down (was the Labels sheet / A3 selection before).
$wsCommands.Activate() | Out-Null
$wsCommands.Range("C66").Select() | Out-Null
